# Region XII_GABALDON.xlsx - add "Status as of July 4, 2025" column with a
# dropdown-validated percentage-bucket value, backed by a new hidden
# "DropdownOptions" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Add the hidden "DropdownOptions" sheet (after Sheet1) with the list
#    of percentage buckets used by the new dropdown column.
# ---------------------------------------------------------------------
$wsOptions = $wb.Worksheets.Add($null, $ws1)
$wsOptions.Name = "DropdownOptions"

$wsOptions.Cells.Item(1, 1).Value = "0% - 10%"
$wsOptions.Cells.Item(2, 1).Value = "11% - 25%"
$wsOptions.Cells.Item(3, 1).Value = "26% - 50%"
$wsOptions.Cells.Item(4, 1).Value = "51% - 75%"
$wsOptions.Cells.Item(5, 1).Value = "76% - 90%"
$wsOptions.Cells.Item(6, 1).Value = "91% - 99%"
$wsOptions.Cells.Item(7, 1).Value = "100%"

$wsOptions.Visible = $false

# ---------------------------------------------------------------------
# 2. Add the new "Status as of July 4, 2025" header in column AU of
#    Sheet1 (no header styling, matching the source edit).
# ---------------------------------------------------------------------
$ws1.Range("AU1").Value = "Status as of July 4, 2025"

# ---------------------------------------------------------------------
# 3. Clean up the stray empty cells that were scattered across rows 2-7
#    (cells that only held an empty inline string).
# ---------------------------------------------------------------------
$emptyRanges = @(
    "N2", "Q2:AA2", "AN2", "AP2:AS2",
    "N3", "Q3:AA3", "AN3", "AP3:AS3",
    "N4", "Q4:AA4", "AN4", "AP4:AS4",
    "N5", "Q5:AA5", "AN5", "AQ5:AS5",
    "AN6:AO6", "AQ6:AS6",
    "N7", "Q7:AA7", "AN7:AO7", "AQ7:AR7"
)
foreach ($rng in $emptyRanges) {
    $ws1.Range($rng).ClearContents()
}

# ---------------------------------------------------------------------
# 4. Apply the dropdown (list) data validation to AU2:AU7, sourced from
#    the DropdownOptions sheet.
# ---------------------------------------------------------------------
$validation = $ws1.Range("AU2:AU7").Validation
$validation.Add(3, 1, 1, 'DropdownOptions!$A$1:$A$7')
$validation.ShowInput = $false
$validation.ShowError = $false
